$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "The industry will use 5G "
$find.Execute() | Out-Null

if ($find.Found) {
    $target = $find.Parent
    $end = $target.End
    $rng = $d.Range($end, $end)
    $rng.Font.Size = 10
    $rng.InsertAfter("to increase the speed & quality of their services")
}
